$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.00084860496271118668
$ws.Range("B1").Value = 0.00084860496271118668
$ws.Range("C1").Value = 0.00084860496271118668
$ws.Range("D1").Value = 0.00084860496271118733
$ws.Range("F1").Value = 0.00084860496271118668
$ws.Range("G1").Value = 0.99915139503728889192
$ws.Range("H1").Value = 0.00084860496271118441
$ws.Range("J1").Value = 0.00084860496271116478

$ws.Range("B2").Value = 0.99915139503728889192
$ws.Range("D2").Value = 0.99915139503728889192
$ws.Range("F2").Value = 0.99915139503728889192
$ws.Range("G2").Value = 0.00084860496271119113
$ws.Range("H2").Value = 0.99915139503728889192
$ws.Range("I2").Value = 0.00084860496271118668
$ws.Range("J2").Value = 0.99915139503728866988
